# Add pseudocode + analysis for the 4th ("Linear Time") algorithm at the
# end of the document, matching the author's commit:
#   "added pseudocode for 4th algorithm / added final algorithm to the doc"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helpers for building WordprocessingML fragments that we can hand to
# Range.InsertXML (the COM equivalent of pasting OOXML at a Range).
# ---------------------------------------------------------------------

function Wrap-PackageXml($bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyInner + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Bold-Run($text, $tabBefore) {
    $tab = if ($tabBefore) { '<w:tab/>' } else { '' }
    $t = if ($text -ne $null) { "<w:t>$text</w:t>" } else { '' }
    return "<w:r><w:rPr><w:b/></w:rPr>$tab$t</w:r>"
}

function Wingdings-EmptyPara {
    return '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings"/></w:rPr></w:pPr></w:p>'
}

function Code-Para($runsXml) {
    return '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/><w:rPr><w:b/></w:rPr></w:pPr>' + $runsXml + '</w:p>'
}

function Code-EmptyPara {
    return '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/><w:rPr><w:b/></w:rPr></w:pPr></w:p>'
}

# Build the pseudocode + analysis block that follows the "Linear Time"
# heading (mirrors the style already used for the other three algorithms
# earlier in the document).
$pseudocode =
    (Code-Para ((Bold-Run "linearTime" $false) + (Bold-Run "(" $false) + (Bold-Run "array" $false) + (Bold-Run ")" $false))) +
    (Code-Para (Bold-Run "loop I = 0 for n" $true)) +
    (Code-Para ((Bold-Run $null $true) + (Bold-Run "tempMax = max(array[i], tempMax + array[i])" $true))) +
    (Code-Para ((Bold-Run $null $true) + (Bold-Run "maxSum = max(maxSum, tempMax)" $true))) +
    (Code-Para (Bold-Run "return maxSum" $true)) +
    (Code-EmptyPara) +
    (Code-Para (Bold-Run "Analysis:" $false)) +
    (Code-EmptyPara)

$analysisText = "This algorithm is not recursive, and you are only looping once for n. Therefore the running time is O(n)."
$analysisPara = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/></w:pPr>' +
    "<w:r><w:tab/><w:t>$analysisText</w:t></w:r>" +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$trailingBlankCodeParas = (Code-EmptyPara) + (Code-EmptyPara)

$finalWingdingsPara = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings"/></w:rPr><w:tab/></w:r></w:p>'

$fullInsert = (Wingdings-EmptyPara) + (Wingdings-EmptyPara) + $pseudocode + $analysisPara +
    $trailingBlankCodeParas + $finalWingdingsPara

# ---------------------------------------------------------------------
# Locate the "Linear Time" heading paragraph (last bullet added for the
# 4th algorithm) and the old trailing blank paragraph that follows it;
# the trailing paragraph currently owns the document's "_GoBack" bookmark,
# which needs to move down to the new analysis paragraph.
# ---------------------------------------------------------------------

$headingRange = $d.Content.Find
$found = $d.Content.Find.Execute("Linear Time")

$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd("`r") -eq "Linear Time") {
        $headingPara = $d.Paragraphs($i)
    }
}

# The "_GoBack" bookmark sits on the paragraph right after the heading;
# remove it now -- it gets re-created on the new analysis paragraph above.
$d.Bookmarks.Item("_GoBack").Delete()

# Insert one blank paragraph right after the heading, then replace its
# contents with the full block of new paragraphs via InsertXML so that
# every paragraph/run keeps the exact formatting used elsewhere in the
# document.
$headRange = $headingPara.Range
$headRange.Collapse(0)
$headRange.InsertParagraphAfter()

$targetPara = $headingPara.Next()
$targetPara.Range.InsertXML((Wrap-PackageXml $fullInsert))
